# Fruta / hortaliza, semanal
# Insert a new weekly data row at row 345 (shifting existing rows 345-366 down to 346-367)
# on the single worksheet of the "Poroto verde" (Hortaliza) workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 345; this shifts rows 345:366 down to 346:367
# and also extends the sheet dimension from A1:R366 to A1:R367 automatically.
$ws.Rows.Item(345).Insert()

# Populate the newly inserted row 345 with the new weekly record.
$ws.Cells.Item(345, 1).Value = 8
$ws.Cells.Item(345, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(345, 3).Value = "Coquimbo"
$ws.Cells.Item(345, 4).Value = 45021
$ws.Cells.Item(345, 5).Value = 4
$ws.Cells.Item(345, 6).Value = 100112031
$ws.Cells.Item(345, 7).Value = "Poroto verde"
$ws.Cells.Item(345, 8).Value = "Magnum"
$ws.Cells.Item(345, 9).Value = "Primera"
$ws.Cells.Item(345, 10).Value = 440
$ws.Cells.Item(345, 11).Value = 22000
$ws.Cells.Item(345, 12).Value = 23000
$ws.Cells.Item(345, 13).Value = 22500
$ws.Cells.Item(345, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(345, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(345, 16).Value = 900
$ws.Cells.Item(345, 17).Value = 25
$ws.Cells.Item(345, 18).Value = "Hortaliza"
